$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column N down into the new column O (row by row,
# since some of the N rows carry different styles from each other), then
# overwrite the values with the 2021 figures.
$dataRows = 4,5,6,7,8,9,10,11,12,13,14,16,17
foreach ($r in $dataRows) {
    $ws.Range("N$r").Copy() | Out-Null
    $ws.Range("O$r").PasteSpecial(-4122) | Out-Null
}

$ws.Range("O4").Value = 2021

$ws.Range("O5").Value = 11.7
$ws.Range("O6").Value = 16.4
$ws.Range("O7").Value = 9.7
$ws.Range("O8").Value = 12.1
$ws.Range("O9").Value = 5.3
$ws.Range("O10").Value = 4.7
$ws.Range("O11").Value = 3.4
$ws.Range("O12").Value = 18.8
$ws.Range("O13").Value = 19.6
$ws.Range("O14").Value = 6.9
$ws.Range("O16").Value = 12.8
$ws.Range("O17").Value = 11

# Update the saved view: scroll back to the top-left and move the
# selection from S18 to R11.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("R11").Select() | Out-Null
